# Apply crypto price/volume update (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "70.352.00"
$ws.Range("E2").Value = "  -3.69%  "

# Row 3
$ws.Range("D3").Value = "3.825.61"
$ws.Range("E3").Value = "  -4.44%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.17%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.16%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "164.77"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.38%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.668"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.92%  "

# Row 8
$ws.Range("E8").Value = "  +0.37%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.741"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.57%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.173"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.69%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.62"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.33%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000317"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.74%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.23"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.69%  "

# Row 14
$ws.Range("D14").Value = "4.445.32"
$ws.Range("E14").Value = "  -4.12%  "

# Row 15
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "3.836.94"
$ws.Range("E15").Value = "  -4.01%  "

# Row 16
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.83"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.89%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.70"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.37%  "

# Row 18
$ws.Range("E18").Value = "  -6.22%  "

# Row 19
$ws.Range("E19").Value = "  -2.32%  "

# Row 20
$ws.Range("D20").Value = "70.364.05"
$ws.Range("E20").Value = "  -3.39%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "433.67"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.69%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.69"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.78%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "93.16"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.40%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.23"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -6.35%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.78"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.04%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.09"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.26%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.91"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -12.80%  "

# Row 28
$ws.Range("E28").Value = "  -0.07%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.35"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.63%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.82"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.57%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.08"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.67%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "13.35"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.90%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "47.88"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.08%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.124"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.74%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "69.26"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.50%  "

# Row 36
$ws.Range("D36").Value = "0.0₃0972"
$ws.Range("E36").Value = "  +10.28%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "634.50"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.87%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.421"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.24%  "

# Row 39
$ws.Range("E39").Value = "  +0.05%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.144"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.82%  "

# Row 41
$ws.Range("E41").Value = "  -0.04%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.24"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.10%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.18"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +18.82%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0466"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.99%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.71"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.04%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.83"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -9.61%  "

# Row 47
$ws.Range("E47").Value = "  -5.05%  "

# Row 48
$ws.Range("E48").Value = "  -14.94%  "

# Row 49
$ws.Range("B49").Value = "ApeXProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.26"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.47%  "

# Row 50
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "2.823.17"
$ws.Range("E50").Value = "  -0.16%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.000270"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.12%  "
